$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 638 (shifts old row 638..659 down to 639..660,
# including style of the D column date cell).
$ws.Rows.Item(638).Insert()

# Populate the newly inserted row 638 with this week's data point
# (same market/product/variety metadata as its neighbours, new
# date/volume/prices for the week).
$ws.Cells.Item(638, 1).Value = 4
$ws.Cells.Item(638, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(638, 3).Value = "Los Lagos"
$ws.Cells.Item(638, 4).Value = 45147
$ws.Cells.Item(638, 5).Value = 10
$ws.Cells.Item(638, 6).Value = "Fruta"
$ws.Cells.Item(638, 7).Value = 100102
$ws.Cells.Item(638, 8).Value = "Cítricos"
$ws.Cells.Item(638, 9).Value = 100102006
$ws.Cells.Item(638, 10).Value = "Pomelo"
$ws.Cells.Item(638, 11).Value = "Start Ruby"
$ws.Cells.Item(638, 12).Value = "Primera"
$ws.Cells.Item(638, 13).Value = 60
$ws.Cells.Item(638, 14).Value = 14000
$ws.Cells.Item(638, 15).Value = 14000
$ws.Cells.Item(638, 16).Value = 14000
$ws.Cells.Item(638, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(638, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(638, 19).Value = 1000
$ws.Cells.Item(638, 20).Value = 14
